# Update gh-pages to output generated at 456a3b4
# Refresh the "want-to-go" counts (column F) and minimum price (column G)
# across the four sheets: 展览 (Exhibition), 演出 (Performance),
# 本地生活 (Local life) and 全部类型 (All types).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 ---
$ws1.Range("F3").Value2 = 974
$ws1.Range("F4").Value2 = 607
$ws1.Range("F5").Value2 = 3027
$ws1.Range("F6").Value2 = 807
$ws1.Range("F9").Value2 = 451
$ws1.Range("F10").Value2 = 660
$ws1.Range("F12").Value2 = 571
$ws1.Range("F13").Value2 = 529
$ws1.Range("F15").Value2 = 1262
$ws1.Range("F16").Value2 = 750
$ws1.Range("F17").Value2 = 32
$ws1.Range("F23").Value2 = 535
$ws1.Range("F25").Value2 = 667
$ws1.Range("F26").Value2 = 667
$ws1.Range("F27").Value2 = 18
$ws1.Range("F30").Value2 = 20
$ws1.Range("F33").Value2 = 229
$ws1.Range("F34").Value2 = 127
$ws1.Range("F35").Value2 = 911
$ws1.Range("F36").Value2 = 4707
$ws1.Range("F37").Value2 = 263
$ws1.Range("F38").Value2 = 43
$ws1.Range("F39").Value2 = 12
$ws1.Range("F40").Value2 = 84

# --- Sheet 2: 演出 ---
# Row 2's price became unavailable ("不可售") instead of a numeric value.
$ws2.Range("G2").Value2 = "不可售"
$ws2.Range("F5").Value2 = 70
$ws2.Range("F8").Value2 = 340
$ws2.Range("F14").Value2 = 171
$ws2.Range("F25").Value2 = 306
$ws2.Range("F27").Value2 = 195
$ws2.Range("F31").Value2 = 31
$ws2.Range("F37").Value2 = 595
$ws2.Range("G37").Value2 = 188

# --- Sheet 3: 本地生活 ---
$ws3.Range("F4").Value2 = 1469
$ws3.Range("F6").Value2 = 277

# --- Sheet 4: 全部类型 ---
$ws4.Range("F3").Value2 = 1469
$ws4.Range("F6").Value2 = 277
$ws4.Range("F7").Value2 = 974
$ws4.Range("F8").Value2 = 607
$ws4.Range("F9").Value2 = 3027
$ws4.Range("F10").Value2 = 807
$ws4.Range("F13").Value2 = 451
$ws4.Range("F14").Value2 = 660
$ws4.Range("F16").Value2 = 70
$ws4.Range("F17").Value2 = 571
$ws4.Range("F19").Value2 = 529
$ws4.Range("F22").Value2 = 1262
$ws4.Range("F23").Value2 = 750
$ws4.Range("F24").Value2 = 32
$ws4.Range("F32").Value2 = 535
$ws4.Range("F35").Value2 = 667
$ws4.Range("F36").Value2 = 667
$ws4.Range("F37").Value2 = 18
$ws4.Range("F39").Value2 = 20
$ws4.Range("F42").Value2 = 306
$ws4.Range("F43").Value2 = 229
$ws4.Range("F45").Value2 = 127
$ws4.Range("F46").Value2 = 911
$ws4.Range("F47").Value2 = 4707
$ws4.Range("F48").Value2 = 263
$ws4.Range("F49").Value2 = 43
$ws4.Range("F50").Value2 = 595
$ws4.Range("G50").Value2 = 188
$ws4.Range("F51").Value2 = 595
$ws4.Range("G51").Value2 = 188
